$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$r = $ws1.Range("A1")
Write-Output $r.Borders.Item(7).LineStyle
Write-Output $r.Interior.ColorIndex
Write-Output $r.Font.Bold
